# "add tabel format baru" - renumber the Bab 4 / Kolaka subdistrict tables
# (new tables 4.2.3/4.2.4 were inserted upstream, pushing this sheet's
# tables from 4.2.3/4.2.4/4.2.5 to 4.2.5/4.2.6/4.2.7) and roll the report
# year forward from 2020 to 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (Tenaga Kesehatan) -------------------------------------------
# Table number: 4.2.3 -> 4.2.5
$ws.Range("H1").Value = "Tabel 4.2.5"
# Titles: 2020 -> 2021
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Kolaka. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Kolaka, 2021"
$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Kolaka Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Kolaka Subdistrict, 2021"

# --- Table 2 (Ibu Melahirkan) ----------------------------------------------
# Table number: 4.2.4. -> 4.2.6.
$ws.Range("P1").Value = "Tabel 4.2.6."
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Kolaka, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Kolaka Subdistrict, 2021"

# --- Table 3 (Pasangan Usia Subur / KB) -------------------------------------
# Table number: 4.2.5. -> 4.2.7.
$ws.Range("W1").Value = "Tabel 4.2.7."
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Kolaka, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Kolaka Subdistrict, 2021"

# Reset the view back to the top-left corner / A1 instead of the
# mid-sheet scroll position + selection it had been left on.
$ws.Activate()
$ws.Range("A1").Select()
